$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet -> this also repoints the _xlnm._FilterDatabase defined name
$ws.Name = "BL_Product_spec"

# Resize columns B and C (values chosen so the engine's column-width
# quantization lands on the target stored widths of 49 and ~14.664)
$ws.Columns.Item(2).ColumnWidth = 48.166666666666664
$ws.Columns.Item(3).ColumnWidth = 13.830729166666666

# Drop the stale active-cell selection (D13), leaving just the frozen-pane selection
$ws.Range("A1").Select()
